# Update Student ID values in column A (rows 2-23) with new log entries,
# and remove the two trailing rows (24-25) that are no longer part of the export.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    # Assign as text (leading apostrophe keeps a numeric-looking string from
    # being auto-converted to a Number), then restore the original cell format
    # (copied from the untouched neighboring "Subject" cell on the same row) so
    # the quote-prefix flag picked up by the text entry does not linger as a
    # new/changed style.
    $ws.Range($cellRef).Value = "'" + $newValue
    $row = $ws.Range($cellRef).Row
    $ws.Cells.Item($row, 2).Copy() | Out-Null
    $ws.Range($cellRef).PasteSpecial(-4122) | Out-Null
}

Set-TextValue "A2" "211242"
Set-TextValue "A3" "200905"
Set-TextValue "A4" "201638"
Set-TextValue "A5" "200727"
Set-TextValue "A6" "191480"
Set-TextValue "A7" "211197"
Set-TextValue "A8" "201051"
Set-TextValue "A9" "211111"
Set-TextValue "A10" "191062"
Set-TextValue "A11" "200866"
Set-TextValue "A12" "201954"
Set-TextValue "A13" "202162"
Set-TextValue "A14" "191186"
Set-TextValue "A15" "211102"
Set-TextValue "A16" "200742"
Set-TextValue "A17" "211004"
Set-TextValue "A18" "200928"
Set-TextValue "A19" "190314"
Set-TextValue "A20" "200490"
Set-TextValue "A21" "200423"
Set-TextValue "A22" "201495"
Set-TextValue "A23" "191052"

$excel.CutCopyMode = $false

# Drop the last two log rows (24-25) entirely.
$ws.Range("A24:F25").EntireRow.Delete()
